$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 1604.037
$ws.Range("I113").Value = 1536.7727
$ws.Range("K113").Value = 1536.7727
$ws.Range("M113").Value = 1717.2273
$ws.Range("H116").Value = 5674.7207
$ws.Range("I116").Value = 7853.316
$ws.Range("J116").Value = 3950
$ws.Range("K116").Value = 7853.316
$ws.Range("L116").Value = 3950
$ws.Range("M116").Value = -4411.316
$ws.Range("N116").Value = -10834
$ws.Range("H125").Value = 45456070
$ws.Range("I125").Value = 90910504
$ws.Range("J125").Value = 1635.6364
$ws.Range("K125").Value = 818194536
$ws.Range("L125").Value = 14720.7276
$ws.Range("M125").Value = -818192076
$ws.Range("N125").Value = -19640.7276
$ws.Range("H132").Value = 120633.234
$ws.Range("I132").Value = 1196.4038
$ws.Range("K132").Value = 3589.2114
$ws.Range("M132").Value = -1059.2114
$ws.Range("H137").Value = 31370.484
$ws.Range("I137").Value = 40970.64
$ws.Range("J137").Value = 1370
$ws.Range("K137").Value = 122911.92
$ws.Range("L137").Value = 4110
$ws.Range("M137").Value = -120361.92
$ws.Range("N137").Value = -9210

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1327.5385
$ws.Range("I61").Value = 1396
$ws.Range("J61").Value = 1173.5
$ws.Range("K61").Value = 1396
$ws.Range("L61").Value = 1173.5
$ws.Range("M61").Value = -1184
$ws.Range("N61").Value = -1597.5
$ws.Range("H74").Value = 52735.562
$ws.Range("I74").Value = 78333.19500000001
$ws.Range("J74").Value = 1540.3077
$ws.Range("K74").Value = 78333.19500000001
$ws.Range("L74").Value = 1540.3077
$ws.Range("M74").Value = -77459.19500000001
$ws.Range("N74").Value = -3288.3077
$ws.Range("H77").Value = 52735.562
$ws.Range("I77").Value = 78333.19500000001
$ws.Range("J77").Value = 1540.3077
$ws.Range("K77").Value = 391665.975
$ws.Range("L77").Value = 7701.538500000001
$ws.Range("M77").Value = -387297.975
$ws.Range("N77").Value = -16437.5385
$ws.Range("H132").Value = 2158882
$ws.Range("I132").Value = 2220622.8
$ws.Range("J132").Value = 1685535.4
$ws.Range("K132").Value = 6661868.399999999
$ws.Range("L132").Value = 5056606.199999999
$ws.Range("M132").Value = -6659338.399999999
$ws.Range("N132").Value = -5061666.199999999
$ws.Range("H133").Value = 22880.143
$ws.Range("J133").Value = 22880.143
$ws.Range("L133").Value = 22880.143
$ws.Range("N133").Value = -27940.143
$ws.Range("H136").Value = 1327.5385
$ws.Range("I136").Value = 1396
$ws.Range("J136").Value = 1173.5
$ws.Range("K136").Value = 4188
$ws.Range("L136").Value = 3520.5
$ws.Range("M136").Value = -1638
$ws.Range("N136").Value = -8620.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 31149.324
$ws.Range("I134").Value = 1064.871
$ws.Range("J134").Value = 186585.67
$ws.Range("K134").Value = 3194.613
$ws.Range("L134").Value = 559757.01
$ws.Range("M134").Value = -659.6130000000003
$ws.Range("N134").Value = -564827.01

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 10330.3125
$ws.Range("J74").Value = 11000
$ws.Range("L74").Value = 11000
$ws.Range("N74").Value = -12748
$ws.Range("H77").Value = 10330.3125
$ws.Range("J77").Value = 11000
$ws.Range("L77").Value = 33000
$ws.Range("N77").Value = -41736
$ws.Range("H99").Value = 2976.6667
$ws.Range("I99").Value = 2950
$ws.Range("J99").Value = 3030
$ws.Range("K99").Value = 2950
$ws.Range("L99").Value = 3030
$ws.Range("M99").Value = -1452
$ws.Range("N99").Value = -6026
$ws.Range("H126").Value = 2976.6667
$ws.Range("I126").Value = 2950
$ws.Range("J126").Value = 3030
$ws.Range("K126").Value = 8850
$ws.Range("L126").Value = 9090
$ws.Range("M126").Value = -6380
$ws.Range("N126").Value = -14030
$ws.Range("H132").Value = 1191.2963
$ws.Range("I132").Value = 886.125
$ws.Range("J132").Value = 3632.6667
$ws.Range("K132").Value = 2658.375
$ws.Range("L132").Value = 10898.0001
$ws.Range("M132").Value = -128.375
$ws.Range("N132").Value = -15958.0001
$ws.Range("H133").Value = 35831.5
$ws.Range("J133").Value = 35831.5
$ws.Range("L133").Value = 35831.5
$ws.Range("N133").Value = -40891.5
$ws.Range("H134").Value = 934.25
$ws.Range("I134").Value = 846.63635
$ws.Range("K134").Value = 2539.90905
$ws.Range("M134").Value = -4.909050000000207

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 8814.157999999999
$ws.Range("I5").Value = 6762.6665
$ws.Range("K5").Value = 20287.9995
$ws.Range("M5").Value = -20175.9995
$ws.Range("H131").Value = 13441708
$ws.Range("J131").Value = 14706763
$ws.Range("L131").Value = 44120289
$ws.Range("N131").Value = -44130369
$ws.Range("H135").Value = 8814.157999999999
$ws.Range("I135").Value = 6762.6665
$ws.Range("K135").Value = 60863.9985
$ws.Range("M135").Value = -58328.9985

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 9800
$ws.Range("J96").Value = 9800
$ws.Range("L96").Value = 9800
$ws.Range("N96").Value = -15292
$ws.Range("H132").Value = 42089.56
$ws.Range("I132").Value = 1610.3077
$ws.Range("J132").Value = 85942.086
$ws.Range("K132").Value = 4830.9231
$ws.Range("L132").Value = 257826.258
$ws.Range("M132").Value = -2300.9231
$ws.Range("N132").Value = -262886.258

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2099.2307
$ws.Range("I68").Value = 1806.25
$ws.Range("J68").Value = 2568
$ws.Range("K68").Value = 1806.25
$ws.Range("L68").Value = 2568
$ws.Range("M68").Value = -1057.25
$ws.Range("N68").Value = -4066
$ws.Range("H71").Value = 2099.2307
$ws.Range("I71").Value = 1806.25
$ws.Range("J71").Value = 2568
$ws.Range("K71").Value = 9031.25
$ws.Range("L71").Value = 12840
$ws.Range("M71").Value = -5287.25
$ws.Range("N71").Value = -20328
$ws.Range("H136").Value = 324473.4
$ws.Range("I136").Value = 501467.9
$ws.Range("J136").Value = 2665.2727
$ws.Range("K136").Value = 1504403.7
$ws.Range("L136").Value = 7995.8181
$ws.Range("M136").Value = -1501853.7
$ws.Range("N136").Value = -13095.8181

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 194812.83
$ws.Range("I136").Value = 2088.7917
$ws.Range("J136").Value = 2507501.2
$ws.Range("K136").Value = 6266.375100000001
$ws.Range("L136").Value = 7522503.600000001
$ws.Range("M136").Value = -3716.375100000001
$ws.Range("N136").Value = -7527603.600000001
$ws.Range("H137").Value = 53286.11
$ws.Range("J137").Value = 53286.11
$ws.Range("L137").Value = 53286.11
$ws.Range("N137").Value = -63486.11
$ws.Range("H138").Value = 50975
$ws.Range("J138").Value = 51300
$ws.Range("L138").Value = 51300
$ws.Range("N138").Value = -61580
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()
